# Update column G ("K" - strikeouts) values for rows 2-14 on Sheet1.
# The data in this save-data workbook is regenerated from an external
# source (regen save_data to use K instead of Strike#); here we simply
# write the freshly computed K values into the existing cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newK = @{
    2  = 2
    3  = 2
    4  = 6
    5  = 4
    6  = 6
    7  = 4
    8  = 4
    9  = 5
    10 = 4
    11 = 5
    12 = 4
    13 = 2
    14 = 3
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
